$wb = $excel.ActiveWorkbook
$infoSheet = $wb.Worksheets.Item("Info")
$cal = $wb.Worksheets.Add($infoSheet)
$cal.Name = "Calibration"

$cal.Range("A1").Value = "Country"
$cal.Range("B1").Value = "Calibration parameter [-]"

$countries = @(
  "Belgium", "Bulgaria", "Czechia", "Denmark", "Germany", "Ireland", "Greece", "Spain",
  "France", "Croatia", "Italy", "Latvia", "Luxembourg", "Hungary", "Netherlands", "Austria",
  "Poland", "Portugal", "Romania", "Slovenia", "Slovakia", "Finland", "Sweden", "United Kingdom",
  "Norway", "Switzerland", "Montenegro", "North Macedonia", "Albania", "Serbia",
  "Bosnia and Herzegovina", "Iceland", "Lithuania", "Estonia"
)
$values = @(
  0.54548975061321947, 0.64565470296552097, 0.83256731656608352, 0.69866925051573792,
  0.7638481298295885, 1.3525028824826275, 0.61877967585985072, 0.72860996139960199,
  0.8465401033692932, 0.85868427395152502, 0.88104732659338392, 0.81364728928513641,
  0.78273477652379519, 0.82030913323635868, 0.91025981698755565, 0.71806010271148657,
  0.7775995676403793, 0.84776127757560216, 1.3843049753180527, 0.65525046911680618,
  0.75473896635266557, 1.2507112866864762, 0.76509948151502349, 0.86455759398530896,
  1.0682268429937749, 1.0432411508595323, 1.3159435219236699, 0.6262913652931289,
  1, 0.90646455597979136, 1, 1, 0.94209056219047782, 0.98423224307957646
)

for ($i = 0; $i -lt $countries.Length; $i++) {
    $row = $i + 2
    $cal.Range("A$row").Value = $countries[$i]
    $cal.Range("B$row").Value = $values[$i]
}

Write-Output "done"
